$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 9; existing rows 9-19 shift down to 10-20.
$ws.Rows(9).Insert()

# Populate the new row 9 with the new weekly record.
$ws.Range("A9").Value = 11
$ws.Range("B9").Value = "Vega Monumental Concepción"
$ws.Range("C9").Value = "Bíobío"
$ws.Range("D9").Value = "12/10/2021"
$ws.Range("E9").Value = 8
$ws.Range("F9").Value = 100112022
$ws.Range("G9").Value = "Arveja Verde"
$ws.Range("H9").Value = "Sin especificar"
$ws.Range("I9").Value = "Primera"
$ws.Range("J9").Value = 110
$ws.Range("K9").Value = 16000
$ws.Range("L9").Value = 17000
$ws.Range("M9").Value = 16545
$ws.Range("N9").Value = "$/saco 25 kilos"
$ws.Range("O9").Value = "Región del Maule"
$ws.Range("P9").Value = 662
$ws.Range("Q9").Value = 25
$ws.Range("R9").Value = "Hortaliza"
